# Add an extra "annotated images" count column (H) plus the ratio of the
# new count to the existing average (I), and an overall average of the
# ratios at the bottom (I12).
#
# The existing data (labels in column A, four repeated counts in B:E and
# their average in F) shifts one column to the right to make room, which
# is exactly what Excel does when a whole column is inserted before A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right: A->B, B->C, ... F->G.
$ws.Columns("A:A").Insert()

# New raw counts for the "annotated images" pass (goes into column H,
# next to the pre-existing - but previously empty - "Buchnera 1" header
# that also shifted from G1 to H1).
$counts = @(146, 170, 162, 126, 79, 88, 126, 152, 61, 227)
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value2 = $counts[$i]
}

# Recreate the row average in column G (was F before the insert), now
# referencing the shifted C:F columns. Filling the whole range at once
# reproduces Excel's "enter formula + fill down" shared-formula grouping.
$ws.Range("G2:G11").Formula = "=AVERAGE(C2,D2,E2,F2)"

# Ratio of the new count (H) to the recomputed average (G). The first row
# is entered on its own (so it stays a standalone formula) and the rest
# are filled down together as a second shared-formula group, matching how
# this was actually built up in Excel.
$ws.Range("I2").Formula = "=H2/G2"
$ws.Range("I3:I11").Formula = "=H3/G3"

# Overall average of the ratio column.
$ws.Range("I12").Formula = "=AVERAGE(I2:I11)"

# Leave the view roughly where the author left it.
$ws.Range("E19").Select()
